# Applies the zh-Hant -> zh-Hans (Traditional -> Simplified) retranslation
# described by the commit diff to the active document.

$d = $word.ActiveDocument

function Global-Replace($doc, $find, $replace) {
    $doc.Content.Find.Execute($find, $false, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
}

function Get-ParaIndexByMarkerNth($doc, $marker, $n) {
    $count = 0
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        if ($p.Range.Text.Contains($marker)) {
            $count = $count + 1
            if ($count -eq $n) {
                return $i
            }
        }
    }
    return -1
}

function Scoped-Replace($doc, $paraIndex, $find, $replace) {
    $p = $doc.Paragraphs.Item($paraIndex)
    $r = $doc.Range($p.Range.Start, $p.Range.End)
    $r.Find.Execute($find, $false, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
}

# ---------------------------------------------------------------------------
# 1) Changes that are identical in every place they occur in the body
#    (safe to Replace-All across the whole document).
# ---------------------------------------------------------------------------

Global-Replace $d "英文" "英语"
Global-Replace $d " / 葡萄牙文 / 法文 / 泰文 / 越南文 / 西班牙文" " / 葡萄牙语 / 法语 / 泰语 / 越南语 / 西班牙语"

Global-Replace $d "簡介" "简要"
Global-Replace $d "發送給在目標國家中已回覆參加但尚未寄送文件的合作夥伴的電子郵件。 將通過 customer.io 發送" "发送给在目标国家中确认参加但尚未向我们提交文件的合作伙伴的电子邮件。 将通过 customer.io 发送"
Global-Replace $d "目標受眾" "目标受众"
Global-Replace $d "尚未提交文件的受邀合作夥伴" "已邀请但尚未提交文件的合作伙伴"

# Colons: "主題行: " / "簡介:" / "目標受眾:" style runs -> full width "：" with no space.
Global-Replace $d ": " "："
Global-Replace $d ":" "："

Global-Replace $d "主題行" "主题行"
Global-Replace $d "[活動名稱]" "[活动名称]"
Global-Replace $d " — 您已提交文件了嗎？  " " — 您是否已提交文件？  "

Global-Replace $d "不要忘記傳送文檔" "不要忘记发送文件"
Global-Replace $d "[合作夥伴姓名]" "[合作伙伴姓名]"

Global-Replace $d "很高興能在即將舉行的 " "很高兴能在即将举行的 "

Global-Replace $d " 之前提供以下文檔：" " 之前提供以下文件："
Global-Replace $d "[插入所需文件清單]" "[插入所需文件列表]"

Global-Replace $d "，郵箱地址為 " "，邮箱地址为 "
Global-Replace $d "[郵箱地址]" "[电子邮件地址]"
Global-Replace $d "[WHATSAPP 號碼]" "[WHATSAPP 号码]"

Global-Replace $d "期待在那裡見到您！" "期待在那里见到您！"

# ---------------------------------------------------------------------------
# 2) Changes that differ between the first ("注册确认") email template and the
#    second ("最佳体验") email template further down in the same document, so
#    they must be scoped to the specific paragraph they live in.
# ---------------------------------------------------------------------------

# " 見到您。 " (plain) vs " 見到您。 '" (stray quote) endings — disambiguate
# the two otherwise-identical paragraphs by which ending they carry.
$idx1 = Get-ParaIndexByMarkerNth $d " 見到您。 ‘" 1
Scoped-Replace $d $idx1 " 見到您。 ‘" " 见到您。 ‘"

$idx2 = Get-ParaIndexByMarkerNth $d " 見到您。 " 1
Scoped-Replace $d $idx2 " 見到您。 " " 见到您。 "

# Section 1 only.
$idx = Get-ParaIndexByMarkerNth $d "為了確認註冊，需要您在 " 1
Scoped-Replace $d $idx "為了確認註冊，需要您在 " "为了确认注册，需要您在 "

# Section 2 only.
$idx = Get-ParaIndexByMarkerNth $d "為了確保您在此次活動中擁有最佳體驗，我們需要您在 " 1
Scoped-Replace $d $idx "為了確保您在此次活動中擁有最佳體驗，我們需要您在 " "为了确保您在此活动中获得最佳体验，我们需要您在 "

# Section 1 only.
$idx = Get-ParaIndexByMarkerNth $d "請將這些文檔的副本傳送給您的區域經理 " 1
Scoped-Replace $d $idx "請將這些文檔的副本傳送給您的區域經理 " "请将这些文件的副本发送给您的区域经理 "

# Section 1 only.
$idx = Get-ParaIndexByMarkerNth $d " (WhatsApp)，以便我們做出必要的安排，包括住宿和交通。" 1
Scoped-Replace $d $idx " (WhatsApp)，以便我們做出必要的安排，包括住宿和交通。" " (WhatsApp)，以便我们做出必要的安排，包括住宿和交通。"

# Section 1 only (has trailing period, distinct from section 2's wording).
$idx = Get-ParaIndexByMarkerNth $d "如有任何疑問，請聯繫您的區域經理。" 1
Scoped-Replace $d $idx "如有任何疑問，請聯繫您的區域經理。" "如有任何疑问，请联系您的区域经理。"

# Section 2 only.
$idx = Get-ParaIndexByMarkerNth $d "請回覆此電子郵件，附上這些文檔的副本，以便我們做出必要的安排，包括住宿和交通。" 1
Scoped-Replace $d $idx "請回覆此電子郵件，附上這些文檔的副本，以便我們做出必要的安排，包括住宿和交通。" "请回复此电子邮件，附上这些文件的副本，以便我们能为您做出必要的安排，包括住宿和交通。"

# Section 2 only.
$idx = Get-ParaIndexByMarkerNth $d "如有任何疑問，請通過 " 1
Scoped-Replace $d $idx "如有任何疑問，請通過 " "如果有任何疑问，请通过 "

$idx = Get-ParaIndexByMarkerNth $d "即時聊天" 1
Scoped-Replace $d $idx "即時聊天" "实时聊天"

$idx = Get-ParaIndexByMarkerNth $d " 與我們聯繫。 " 1
Scoped-Replace $d $idx " 與我們聯繫。 " " 联系我们。 "

# Section 2 only (trailing space, no period — distinct from section 1's wording).
$idx = Get-ParaIndexByMarkerNth $d "如有任何疑問，請聯繫您的區域經理 " 1
Scoped-Replace $d $idx "如有任何疑問，請聯繫您的區域經理 " "如果您有任何问题，请联系您的区域经理 "

# NOTE: "，郵箱地址為 " / "[郵箱地址]" / "[WHATSAPP 號碼]" get the *same*
# replacement text in both section 1 and section 2, so the earlier
# Global-Replace calls already cover both occurrences correctly.

# ---------------------------------------------------------------------------
# 3) Comment text.
# ---------------------------------------------------------------------------

$comments = $d.Comments
for ($i = 1; $i -le $comments.Count; $i++) {
    $c = $comments.Item($i)
    $c.Range.Find.Execute("選擇任一", $false, $false, $false, $false, $false, $true, 1, $false, "选择任一", 2)
}
